$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "24.636.38"
$ws.Range("E2").Value = "  +0.24%  "

# Row 3
$ws.Range("D3").Value = "1.696.26"
$ws.Range("E3").Value = "  +0.07%  "

# Row 4
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").Value = "'315.75"
$ws.Range("E5").Value = "  -0.23%  "

# Row 7
$ws.Range("D7").Value = "'0.3923"
$ws.Range("E7").Value = "  -0.38%  "

# Row 8
$ws.Range("D8").Value = "'0.4048"
$ws.Range("E8").Value = "  +0.64%  "

# Row 9
$ws.Range("D9").Value = "'1.522"
$ws.Range("E9").Value = "  +0.04%  "

# Row 10
$ws.Range("D10").Value = "'1.004"
$ws.Range("E10").Value = "  +0.26%  "

# Row 11
$ws.Range("D11").Value = "'53.09"
$ws.Range("E11").Value = "  -1.35%  "

# Row 12
$ws.Range("D12").Value = "'0.08836"
$ws.Range("E12").Value = "  +0.74%  "

# Row 13
$ws.Range("D13").Value = "'7.429"
$ws.Range("E13").Value = "  +2.82%  "

# Row 14
$ws.Range("D14").Value = "'23.68"
$ws.Range("E14").Value = "  +1.95%  "

# Row 15
$ws.Range("D15").Value = "'8.107"
$ws.Range("E15").Value = "  +6.83%  "

# Row 16
$ws.Range("D16").Value = "'0.00001318"
$ws.Range("E16").Value = "  -0.54%  "

# Row 17
$ws.Range("D17").Value = "1.701.40"
$ws.Range("E17").Value = "  +0.03%  "

# Row 18
$ws.Range("D18").Value = "'99.34"
$ws.Range("E18").Value = "  -0.80%  "

# Row 19
$ws.Range("D19").Value = "'0.07021"
$ws.Range("E19").Value = "  -0.67%  "

# Row 20
$ws.Range("D20").Value = "'19.73"
$ws.Range("E20").Value = "  +0.35%  "

# Row 22
$ws.Range("D22").Value = "'1.005"
$ws.Range("E22").Value = "  +0.41%  "

# Row 23
$ws.Range("D23").Value = "'14.69"
$ws.Range("E23").Value = "  +4.55%  "

# Row 24
$ws.Range("D24").Value = "24.628.38"
$ws.Range("E24").Value = "  +0.26%  "

# Row 25
$ws.Range("D25").Value = "'3.136"
$ws.Range("E25").Value = "  +3.40%  "

# Row 26
$ws.Range("E26").Value = "  +1.50%  "

# Row 27
$ws.Range("D27").Value = "'22.60"
$ws.Range("E27").Value = "  +0.97%  "

# Row 28
$ws.Range("E28").Value = "  +2.17%  "

# Row 29
$ws.Range("D29").Value = "'8.795"
$ws.Range("E29").Value = "  +18.34%  "

# Row 30
$ws.Range("D30").Value = "'135.59"
$ws.Range("E30").Value = "  +0.81%  "

# Row 31
$ws.Range("D31").Value = "'5.126"
$ws.Range("E31").Value = "  -1.67%  "

# Row 32
$ws.Range("D32").Value = "'0.09006"
$ws.Range("E32").Value = "  +5.66%  "

# Row 33
$ws.Range("D33").Value = "'7.618"
$ws.Range("E33").Value = "  +5.02%  "

# Row 34
$ws.Range("D34").Value = "'1.066"
$ws.Range("E34").Value = "  -3.65%  "

# Row 35
$ws.Range("D35").Value = "'1.958"
$ws.Range("E35").Value = "  +0.37%  "

# Row 36
$ws.Range("D36").Value = "'11.02"
$ws.Range("E36").Value = "  -2.99%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02925"
$ws.Range("E37").Value = "  +6.16%  "

# Row 38
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "'0.2749"
$ws.Range("E38").Value = "  +0.57%  "

# Row 39
$ws.Range("D39").Value = "'14.42"
$ws.Range("E39").Value = "  -0.40%  "

# Row 40
$ws.Range("D40").Value = "'0.09162"
$ws.Range("E40").Value = "  +1.13%  "

# Row 41
$ws.Range("D41").Value = "'1.455"
$ws.Range("E41").Value = "  -0.42%  "

# Row 42
$ws.Range("D42").Value = "'0.7650"
$ws.Range("E42").Value = "  -0.76%  "

# Row 43
$ws.Range("D43").Value = "'16.08"
$ws.Range("E43").Value = "  +3.93%  "

# Row 44
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").Value = "'2.589"
$ws.Range("E44").Value = "  +2.23%  "

# Row 45
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.7168"
$ws.Range("E45").Value = "  -0.40%  "

# Row 46
$ws.Range("D46").Value = "'4.205"
$ws.Range("E46").Value = "  -0.19%  "

# Row 47
$ws.Range("E47").Value = "  +0.01%  "

# Row 48
$ws.Range("D48").Value = "'1.329"
$ws.Range("E48").Value = "  -1.87%  "

# Row 49
$ws.Range("D49").Value = "'139.76"
$ws.Range("E49").Value = "  -1.05%  "

# Row 50
$ws.Range("D50").Value = "'0.07970"
$ws.Range("E50").Value = "  -0.66%  "

# Row 51
$ws.Range("D51").Value = "'90.29"
$ws.Range("E51").Value = "  +2.01%  "
